$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the pdfs list cells: remove the space after the comma separators
# (Gallery Nav Images -> full color updated PDFs)
$ws.Range("M2").Value = "2016.pdf,2017.pdf,2018.pdf"
$ws.Range("M3").Value = "2017-2018.pdf,2018-2019.pdf"

# Update the active view / selection to match the author's saved state
$ws.Range("M4").Select()
$excel.ActiveWindow.ScrollColumn = 11
